$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 878.5397  # H129: 887.02856 -> 878.5397
$ws.Cells.Item(129, 9).Value = 592.5  # I129: 626.6667 -> 592.5
$ws.Cells.Item(129, 10).Value = 897.9322  # J129: 898.6866 -> 897.9322
$ws.Cells.Item(129, 11).Value = 1777.5  # K129: 1880.0001 -> 1777.5
$ws.Cells.Item(129, 12).Value = 2693.7966  # L129: 2696.0598 -> 2693.7966
$ws.Cells.Item(129, 13).Value = 3222.5  # M129: 3119.9999 -> 3222.5
$ws.Cells.Item(129, 14).Value = -12693.7966  # N129: -12696.0598 -> -12693.7966

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4513  # H138: 4421.7896 -> 4513
$ws.Cells.Item(138, 9).Value = 10500  # I138: 10600 -> 10500
$ws.Cells.Item(138, 10).Value = 3764.625  # J138: 3694.9412 -> 3764.625
$ws.Cells.Item(138, 11).Value = 31500  # K138: 31800 -> 31500
$ws.Cells.Item(138, 12).Value = 11293.875  # L138: 11084.8236 -> 11293.875
$ws.Cells.Item(138, 13).Value = -26360  # M138: -26660 -> -26360
$ws.Cells.Item(138, 14).Value = -21573.875  # N138: -21364.8236 -> -21573.875

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2654.2222  # H2: 2254.4 -> 2654.2222
$ws.Cells.Item(2, 9).Value = 1966.6666  # I2: 1218.7142 -> 1966.6666
$ws.Cells.Item(2, 10).Value = 2998  # J2: 4671 -> 2998
$ws.Cells.Item(2, 11).Value = 1966.6666  # K2: 1218.7142 -> 1966.6666
$ws.Cells.Item(2, 12).Value = 2998  # L2: 4671 -> 2998
$ws.Cells.Item(2, 13).Value = -1853.6666  # M2: -1105.7142 -> -1853.6666
$ws.Cells.Item(2, 14).Value = -3224  # N2: -4897 -> -3224

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46, 8).Value = 7287.75  # H46: 8430 -> 7287.75
$ws.Cells.Item(46, 9).Value = 0  # I46: 8138 -> 0
$ws.Cells.Item(46, 10).Value = 7287.75  # J46: 8576 -> 7287.75
$ws.Cells.Item(46, 11).Value = 0  # K46: 8138 -> 0
$ws.Cells.Item(46, 12).Value = 7287.75  # L46: 8576 -> 7287.75
$ws.Cells.Item(46, 13).ClearContents()  # M46: -7819 -> (removed)
$ws.Cells.Item(46, 14).Value = -7925.75  # N46: -9214 -> -7925.75

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 641.1142599999999  # H74: 651.0833 -> 641.1142599999999
$ws.Cells.Item(74, 10).Value = 1361.8182  # J74: 1331.6666 -> 1361.8182
$ws.Cells.Item(74, 12).Value = 1361.8182  # L74: 1331.6666 -> 1361.8182
$ws.Cells.Item(74, 14).Value = -3109.8182  # N74: -3079.6666 -> -3109.8182

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 641.1142599999999  # H77: 651.0833 -> 641.1142599999999
$ws.Cells.Item(77, 10).Value = 1361.8182  # J77: 1331.6666 -> 1361.8182
$ws.Cells.Item(77, 12).Value = 6809.090999999999  # L77: 6658.333000000001 -> 6809.090999999999
$ws.Cells.Item(77, 14).Value = -15545.091  # N77: -15394.333 -> -15545.091

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2654.2222  # H116: 2254.4 -> 2654.2222
$ws.Cells.Item(116, 9).Value = 1966.6666  # I116: 1218.7142 -> 1966.6666
$ws.Cells.Item(116, 10).Value = 2998  # J116: 4671 -> 2998
$ws.Cells.Item(116, 11).Value = 1966.6666  # K116: 1218.7142 -> 1966.6666
$ws.Cells.Item(116, 12).Value = 2998  # L116: 4671 -> 2998
$ws.Cells.Item(116, 13).Value = 327.3334  # M116: 1075.2858 -> 327.3334
$ws.Cells.Item(116, 14).Value = -7586  # N116: -9259 -> -7586

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2219.5  # H122: 1580.1666 -> 2219.5
$ws.Cells.Item(122, 9).Value = 2159.4614  # I122: 1518.4348 -> 2159.4614
$ws.Cells.Item(122, 11).Value = 6478.3842  # K122: 4555.3044 -> 6478.3842
$ws.Cells.Item(122, 13).Value = -4028.3842  # M122: -2105.3044 -> -4028.3842

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 14392.23  # H132: 15974.771 -> 14392.23
$ws.Cells.Item(132, 9).Value = 1118.25  # I132: 1193.24 -> 1118.25
$ws.Cells.Item(132, 10).Value = 48180.547  # J132: 52928.6 -> 48180.547
$ws.Cells.Item(132, 11).Value = 3354.75  # K132: 3579.72 -> 3354.75
$ws.Cells.Item(132, 12).Value = 144541.641  # L132: 158785.8 -> 144541.641
$ws.Cells.Item(132, 13).Value = -824.75  # M132: -1049.72 -> -824.75
$ws.Cells.Item(132, 14).Value = -149601.641  # N132: -163845.8 -> -149601.641

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2654.2222  # H3: 2254.4 -> 2654.2222
$ws.Cells.Item(3, 9).Value = 1966.6666  # I3: 1218.7142 -> 1966.6666
$ws.Cells.Item(3, 10).Value = 2998  # J3: 4671 -> 2998
$ws.Cells.Item(3, 11).Value = 1966.6666  # K3: 1218.7142 -> 1966.6666
$ws.Cells.Item(3, 12).Value = 2998  # L3: 4671 -> 2998
$ws.Cells.Item(3, 13).Value = -1852.6666  # M3: -1104.7142 -> -1852.6666
$ws.Cells.Item(3, 14).Value = -3226  # N3: -4899 -> -3226

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2501726.8  # H105: 2382630.2 -> 2501726.8
$ws.Cells.Item(105, 9).Value = 1675.8182  # I105: 1723.4 -> 1675.8182
$ws.Cells.Item(105, 10).Value = 5557344.5  # J105: 4547091 -> 5557344.5
$ws.Cells.Item(105, 11).Value = 1675.8182  # K105: 1723.4 -> 1675.8182
$ws.Cells.Item(105, 12).Value = 5557344.5  # L105: 4547091 -> 5557344.5
$ws.Cells.Item(105, 13).Value = 71.18180000000007  # M105: 23.59999999999991 -> 71.18180000000007
$ws.Cells.Item(105, 14).Value = -5560838.5  # N105: -4550585 -> -5560838.5

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1132.6666  # H16: 1157.6666 -> 1132.6666
$ws.Cells.Item(16, 10).Value = 1399.75  # J16: 1474.75 -> 1399.75
$ws.Cells.Item(16, 12).Value = 1399.75  # L16: 1474.75 -> 1399.75
$ws.Cells.Item(16, 14).Value = -1973.75  # N16: -2048.75 -> -1973.75

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 9).Value = 10868.4  # I31: 10533.936 -> 10868.4
$ws.Cells.Item(31, 10).Value = 3666.4119  # J31: 3864.3125 -> 3666.4119
$ws.Cells.Item(31, 11).Value = 10868.4  # K31: 10533.936 -> 10868.4
$ws.Cells.Item(31, 12).Value = 3666.4119  # L31: 3864.3125 -> 3666.4119
$ws.Cells.Item(31, 13).Value = -10573.4  # M31: -10238.936 -> -10573.4
$ws.Cells.Item(31, 14).Value = -4256.4119  # N31: -4454.3125 -> -4256.4119

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 9).Value = 10868.4  # I34: 10533.936 -> 10868.4
$ws.Cells.Item(34, 10).Value = 3666.4119  # J34: 3864.3125 -> 3666.4119
$ws.Cells.Item(34, 11).Value = 10868.4  # K34: 10533.936 -> 10868.4
$ws.Cells.Item(34, 12).Value = 3666.4119  # L34: 3864.3125 -> 3666.4119
$ws.Cells.Item(34, 13).Value = -10666.4  # M34: -10331.936 -> -10666.4
$ws.Cells.Item(34, 14).Value = -4070.4119  # N34: -4268.3125 -> -4070.4119

# CRP row 38
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(38, 8).Value = 6500  # H38: 0 -> 6500
$ws.Cells.Item(38, 9).Value = 6500  # I38: 0 -> 6500
$ws.Cells.Item(38, 11).Value = 6500  # K38: 0 -> 6500
$ws.Cells.Item(38, 13).Value = -6123  # M38: None -> -6123

# CRP row 46
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(46, 8).Value = 6500  # H46: 0 -> 6500
$ws.Cells.Item(46, 9).Value = 6500  # I46: 0 -> 6500
$ws.Cells.Item(46, 11).Value = 6500  # K46: 0 -> 6500
$ws.Cells.Item(46, 13).Value = -6289  # M46: None -> -6289

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1132.6666  # H113: 1157.6666 -> 1132.6666
$ws.Cells.Item(113, 10).Value = 1399.75  # J113: 1474.75 -> 1399.75
$ws.Cells.Item(113, 12).Value = 1399.75  # L113: 1474.75 -> 1399.75
$ws.Cells.Item(113, 14).Value = -5739.75  # N113: -5814.75 -> -5739.75

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 17468.637  # H132: 3229.524 -> 17468.637
$ws.Cells.Item(132, 9).Value = 22576.043  # I132: 1062.6154 -> 22576.043
$ws.Cells.Item(132, 10).Value = 5721.6  # J132: 6750.75 -> 5721.6
$ws.Cells.Item(132, 11).Value = 67728.129  # K132: 3187.8462 -> 67728.129
$ws.Cells.Item(132, 12).Value = 17164.8  # L132: 20252.25 -> 17164.8
$ws.Cells.Item(132, 13).Value = -65198.129  # M132: -657.8462 -> -65198.129
$ws.Cells.Item(132, 14).Value = -22224.8  # N132: -25312.25 -> -22224.8

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 918.19446  # H68: 927.2857 -> 918.19446
$ws.Cells.Item(68, 10).Value = 1214.1666  # J68: 1250.2941 -> 1214.1666
$ws.Cells.Item(68, 12).Value = 3642.4998  # L68: 3750.8823 -> 3642.4998
$ws.Cells.Item(68, 14).Value = -5264.4998  # N68: -5372.8823 -> -5264.4998

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 918.19446  # H71: 927.2857 -> 918.19446
$ws.Cells.Item(71, 10).Value = 1214.1666  # J71: 1250.2941 -> 1214.1666
$ws.Cells.Item(71, 12).Value = 10927.4994  # L71: 11252.6469 -> 10927.4994
$ws.Cells.Item(71, 14).Value = -19039.4994  # N71: -19364.6469 -> -19039.4994

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1750  # H75: 3206.077 -> 1750
$ws.Cells.Item(75, 9).Value = 750  # I75: 1124.75 -> 750
$ws.Cells.Item(75, 10).Value = 2750  # J75: 4131.1113 -> 2750
$ws.Cells.Item(75, 11).Value = 2250  # K75: 3374.25 -> 2250
$ws.Cells.Item(75, 12).Value = 8250  # L75: 12393.3339 -> 8250
$ws.Cells.Item(75, 13).Value = -1252  # M75: -2376.25 -> -1252
$ws.Cells.Item(75, 14).Value = -10246  # N75: -14389.3339 -> -10246

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 1750  # H78: 3206.077 -> 1750
$ws.Cells.Item(78, 9).Value = 750  # I78: 1124.75 -> 750
$ws.Cells.Item(78, 10).Value = 2750  # J78: 4131.1113 -> 2750
$ws.Cells.Item(78, 11).Value = 6750  # K78: 10122.75 -> 6750
$ws.Cells.Item(78, 12).Value = 24750  # L78: 37180.00169999999 -> 24750
$ws.Cells.Item(78, 13).Value = -1758  # M78: -5130.75 -> -1758
$ws.Cells.Item(78, 14).Value = -34734  # N78: -47164.00169999999 -> -34734

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 4106.1816  # H81: 4363.1113 -> 4106.1816
$ws.Cells.Item(81, 10).Value = 4883.8887  # J81: 5436.4287 -> 4883.8887
$ws.Cells.Item(81, 12).Value = 14651.6661  # L81: 16309.2861 -> 14651.6661
$ws.Cells.Item(81, 14).Value = -16897.6661  # N81: -18555.2861 -> -16897.6661

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(84, 8).Value = 4106.1816  # H84: 4363.1113 -> 4106.1816
$ws.Cells.Item(84, 10).Value = 4883.8887  # J84: 5436.4287 -> 4883.8887
$ws.Cells.Item(84, 12).Value = 43954.99830000001  # L84: 48927.85830000001 -> 43954.99830000001
$ws.Cells.Item(84, 14).Value = -55186.99830000001  # N84: -60159.85830000001 -> -55186.99830000001

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 949.5  # H86: 899 -> 949.5
$ws.Cells.Item(86, 10).Value = 1000  # J86: 0 -> 1000
$ws.Cells.Item(86, 12).Value = 3000  # L86: 0 -> 3000
$ws.Cells.Item(86, 14).Value = -5372  # N86: None -> -5372

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 9886.666999999999  # H87: 11085 -> 9886.666999999999
$ws.Cells.Item(87, 9).Value = 580  # I87: 636 -> 580
$ws.Cells.Item(87, 11).Value = 1740  # K87: 1908 -> 1740
$ws.Cells.Item(87, 13).Value = -492  # M87: -660 -> -492

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 949.5  # H89: 899 -> 949.5
$ws.Cells.Item(89, 10).Value = 1000  # J89: 0 -> 1000
$ws.Cells.Item(89, 12).Value = 9000  # L89: 0 -> 9000
$ws.Cells.Item(89, 14).Value = -20856  # N89: None -> -20856

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 9886.666999999999  # H90: 11085 -> 9886.666999999999
$ws.Cells.Item(90, 9).Value = 580  # I90: 636 -> 580
$ws.Cells.Item(90, 11).Value = 5220  # K90: 5724 -> 5220
$ws.Cells.Item(90, 13).Value = 1020  # M90: 516 -> 1020

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 805.03  # H131: 806.11224 -> 805.03
$ws.Cells.Item(131, 10).Value = 822.3196  # J131: 823.8 -> 822.3196
$ws.Cells.Item(131, 12).Value = 2466.9588  # L131: 2471.4 -> 2466.9588
$ws.Cells.Item(131, 14).Value = -12546.9588  # N131: -12551.4 -> -12546.9588

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5212396.5  # H70: 5212107 -> 5212396.5
$ws.Cells.Item(70, 9).Value = 4108.2856  # I70: 3785.125 -> 4108.2856
$ws.Cells.Item(70, 10).Value = 12504000  # J70: 15628750 -> 12504000
$ws.Cells.Item(70, 11).Value = 4108.2856  # K70: 3785.125 -> 4108.2856
$ws.Cells.Item(70, 12).Value = 12504000  # L70: 15628750 -> 12504000
$ws.Cells.Item(70, 13).Value = -3838.2856  # M70: -3515.125 -> -3838.2856
$ws.Cells.Item(70, 14).Value = -12504540  # N70: -15629290 -> -12504540

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5212396.5  # H73: 5212107 -> 5212396.5
$ws.Cells.Item(73, 9).Value = 4108.2856  # I73: 3785.125 -> 4108.2856
$ws.Cells.Item(73, 10).Value = 12504000  # J73: 15628750 -> 12504000
$ws.Cells.Item(73, 11).Value = 4108.2856  # K73: 3785.125 -> 4108.2856
$ws.Cells.Item(73, 12).Value = 12504000  # L73: 15628750 -> 12504000
$ws.Cells.Item(73, 13).Value = -3172.2856  # M73: -2849.125 -> -3172.2856
$ws.Cells.Item(73, 14).Value = -12505872  # N73: -15630622 -> -12505872

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2055.6843  # H97: 2139.7222 -> 2055.6843
$ws.Cells.Item(97, 9).Value = 972  # I97: 1015.6923 -> 972
$ws.Cells.Item(97, 10).Value = 4403.6665  # J97: 5062.2 -> 4403.6665
$ws.Cells.Item(97, 11).Value = 972  # K97: 1015.6923 -> 972
$ws.Cells.Item(97, 12).Value = 4403.6665  # L97: 5062.2 -> 4403.6665
$ws.Cells.Item(97, 13).Value = -476  # M97: -519.6923 -> -476
$ws.Cells.Item(97, 14).Value = -5395.6665  # N97: -6054.2 -> -5395.6665

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 19376.871  # H132: 19998.7 -> 19376.871
$ws.Cells.Item(132, 9).Value = 3670.4  # I132: 3825.5789 -> 3670.4
$ws.Cells.Item(132, 11).Value = 11011.2  # K132: 11476.7367 -> 11011.2
$ws.Cells.Item(132, 13).Value = -8481.200000000001  # M132: -8946.736699999999 -> -8481.200000000001

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 10).Value = 500000  # J2: 0 -> 500000
$ws.Cells.Item(2, 12).Value = 500000  # L2: 0 -> 500000
$ws.Cells.Item(2, 14).Value = -500224  # N2: None -> -500224

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 805608.8  # H132: 3381.25 -> 805608.8
$ws.Cells.Item(132, 9).Value = 928625.7  # I132: 2508.6667 -> 928625.7
$ws.Cells.Item(132, 11).Value = 2785877.1  # K132: 7526.000100000001 -> 2785877.1
$ws.Cells.Item(132, 13).Value = -2783347.1  # M132: -4996.000100000001 -> -2783347.1

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1285.6666  # H81: 2042.8572 -> 1285.6666
$ws.Cells.Item(81, 9).Value = 1242.8  # I81: 2060 -> 1242.8
$ws.Cells.Item(81, 10).Value = 1500  # J81: 2000 -> 1500
$ws.Cells.Item(81, 11).Value = 2485.6  # K81: 4120 -> 2485.6
$ws.Cells.Item(81, 12).Value = 3000  # L81: 4000 -> 3000
$ws.Cells.Item(81, 13).Value = -1424.6  # M81: -3059 -> -1424.6
$ws.Cells.Item(81, 14).Value = -5122  # N81: -6122 -> -5122

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 1285.6666  # H84: 2042.8572 -> 1285.6666
$ws.Cells.Item(84, 9).Value = 1242.8  # I84: 2060 -> 1242.8
$ws.Cells.Item(84, 10).Value = 1500  # J84: 2000 -> 1500
$ws.Cells.Item(84, 11).Value = 12428  # K84: 20600 -> 12428
$ws.Cells.Item(84, 12).Value = 15000  # L84: 20000 -> 15000
$ws.Cells.Item(84, 13).Value = -7124  # M84: -15296 -> -7124
$ws.Cells.Item(84, 14).Value = -25608  # N84: -30608 -> -25608

# WVR row 121
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(121, 8).Value = 0  # H121: 45750 -> 0
$ws.Cells.Item(121, 10).Value = 0  # J121: 45750 -> 0
$ws.Cells.Item(121, 12).Value = 0  # L121: 45750 -> 0
$ws.Cells.Item(121, 14).ClearContents()  # N121: -49244 -> (removed)
